$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out old E5:F5 values
$ws.Range("E5:F5").ClearContents()

# Fill C5:D9 with zeros
$ws.Range("C5:D9").Value = 0.0
